# The generic "None" label was ambiguous across the equipment / drill / kick
# rows, so the wording is being made more specific:
#   - the "equipment" parameter's "None" option becomes "No equipment"
#     (this label is shown both in the header row and in the value cell)
#   - the "drill" and "kick" parameters get their own distinct "No drill" /
#     "No kick" wording instead of sharing the generic "None" text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (column group label) and value cell for "equipment"
$ws.Range("J3").Value  = "No equipment"
$ws.Range("B11").Value = "No equipment"

# Value cells for "drill" and "kick"
$ws.Range("B14").Value = "No drill"
$ws.Range("B16").Value = "No kick"
